$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated experiment results for random_forest, lsboost, neural_network
# (old_model row, row 5, is unchanged)

# random_forest - row 2
$ws.Range("B2").Value = 57.080564413808737
$ws.Range("C2").Value = 0.43303145802823562
$ws.Range("D2").Value = 40.29733989091374
$ws.Range("E2").Value = 0.48535495991138794
$ws.Range("F2").Value = 0.69667421361163351
$ws.Range("G2").Value = 0.67876841622979556
$ws.Range("H2").Value = 0.514645040088612
$ws.Range("I2").Value = 0.71950071203295229

# lsboost - row 3
$ws.Range("B3").Value = 57.507675829328143
$ws.Range("C3").Value = 0.43627166213101998
$ws.Range("D3").Value = 39.87109185972411
$ws.Range("E3").Value = 0.4926455752264578
$ws.Range("F3").Value = 0.70188715277205194
$ws.Range("G3").Value = 0.67158869414801881
$ws.Range("H3").Value = 0.5073544247735422
$ws.Range("I3").Value = 0.7125290597160111

# neural_network - row 4
$ws.Range("B4").Value = 62.510332793576197
$ws.Range("C4").Value = 0.47422342139427237
$ws.Range("D4").Value = 44.817700701314081
$ws.Range("E4").Value = 0.58208522606355118
$ws.Range("F4").Value = 0.76294510029460916
$ws.Range("G4").Value = 0.75490937631223787
$ws.Range("H4").Value = 0.41791477393644882
$ws.Range("I4").Value = 0.65313901048550294
